$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that used to sit after "muestra mal
#    la demora en la tabla" (it gets re-created later, at the new cursor
#    position inside the "logueo" text).
# ---------------------------------------------------------------------------
try {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete() | Out-Null
} catch {
    # not present - nothing to do
}

# ---------------------------------------------------------------------------
# 2) Insert the two new bullet items (plus their surrounding blank list
#    paragraphs) right after "no muestra la imagen del pedido cuando los
#    lista". We collapse a duplicate of the *existing* blank paragraph that
#    already follows it to its Start and replace that single paragraph with
#    [blank][spinner item][blank][alta producto item][blank] - the trailing
#    blank paragraph re-creates the one we just overwrote, so nothing is
#    lost.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("no muestra la imagen del pedido cuando los lista", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara = $anchor.Paragraphs(1)
$blankPara = $anchorPara.Next()

$insertionPoint = $blankPara.Range.Duplicate()
$insertionPoint.Collapse(1)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newParagraphsXml = (
    '<w:p ' + $ns + '>' +
      '<w:pPr>' +
        '<w:pStyle w:val="Prrafodelista"/>' +
        '<w:rPr><w:u w:val="single"/></w:rPr>' +
      '</w:pPr>' +
    '</w:p>' +
    '<w:p ' + $ns + '>' +
      '<w:pPr>' +
        '<w:pStyle w:val="Prrafodelista"/>' +
        '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
        '<w:rPr><w:u w:val="single"/></w:rPr>' +
      '</w:pPr>' +
      '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">poner el </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>spinner</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> en la espera del </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>logueo</w:t></w:r>' +
      '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
      '<w:bookmarkEnd w:id="0"/>' +
      '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>' +
    '<w:p ' + $ns + '>' +
      '<w:pPr>' +
        '<w:pStyle w:val="Prrafodelista"/>' +
        '<w:rPr><w:u w:val="single"/></w:rPr>' +
      '</w:pPr>' +
    '</w:p>' +
    '<w:p ' + $ns + '>' +
      '<w:pPr>' +
        '<w:pStyle w:val="Prrafodelista"/>' +
        '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
        '<w:rPr><w:color w:val="FF0000"/><w:u w:val="single"/></w:rPr>' +
      '</w:pPr>' +
      '<w:r><w:rPr><w:color w:val="FF0000"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">el alta de producto </w:t></w:r>' +
      '<w:r><w:rPr><w:color w:val="FF0000"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">manda el precio como </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:rPr><w:color w:val="FF0000"/><w:u w:val="single"/></w:rPr><w:t>string</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:rPr><w:color w:val="FF0000"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> y tiene que ser number</w:t></w:r>' +
    '</w:p>' +
    '<w:p ' + $ns + '>' +
      '<w:pPr>' +
        '<w:pStyle w:val="Prrafodelista"/>' +
        '<w:rPr><w:u w:val="single"/></w:rPr>' +
      '</w:pPr>' +
    '</w:p>'
)

$insertionPoint.InsertXML($newParagraphsXml) | Out-Null

# ---------------------------------------------------------------------------
# 3) Mark the "SUGERENCIAS:" paragraph as starting a rendered page (adds a
#    <w:lastRenderedPageBreak/> before its text run).
# ---------------------------------------------------------------------------
$sugRange = $d.Content
$sugRange.Find.Execute("SUGERENCIAS:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sugPara = $sugRange.Paragraphs(1)

$sugInsertionPoint = $sugPara.Range.Duplicate()
$sugInsertionPoint.Collapse(1)
$sugXml = '<w:p ' + $ns + '><w:r><w:lastRenderedPageBreak/><w:t>SUGERENCIAS:</w:t></w:r></w:p>'
$sugInsertionPoint.InsertXML($sugXml) | Out-Null
